$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at AH (34), pushing chlorophyll..water_content's
# columns (old AH..CW) one column to the right (new AI..CX). This also
# correctly bumps each row's "spans" attribute.
$ws.Columns("AH").Insert()

$newColText = "Name of source institute and unique culture identifier. See the description for the proper format and list of allowed institutes, http://www.insdc.org/controlled-vocabulary-culturecollection-qualifier"

# Comments are NOT moved by the column insert above, so move each existing
# comment from its old column to the (now shifted) column one to the right,
# working from the rightmost column back to AH so we never clobber a
# not-yet-copied comment.
for ($col = 101; $col -ge 34; $col--) {
    $srcCell = $ws.Cells.Item(15, $col)
    $dstCell = $ws.Cells.Item(15, $col + 1)
    $text = $srcCell.Comment.Text()
    if ($dstCell.Comment -eq $null) {
        $newCmt = $dstCell.AddComment($text)
        $newCmt.Author = ""
    } else {
        $dstCell.Comment.Text($text)
    }
}

# AH15's original comment ("density of sample") has now been copied over to
# AI15, so overwrite AH15's comment with the new culture_collection text.
$ws.Cells.Item(15, 34).Comment.Text($newColText)

# Give the new AH15 header cell its label (adds the new shared string).
$ws.Cells.Item(15, 34).Value = "culture_collection"
